# Update Name of Algo
# Applies updated numeric results for RandomForest imputation output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -6.994099999999994
$ws.Range("B12").Value = 5.190499999999998
$ws.Range("D14").Value = -7.415399999999998
$ws.Range("D26").Value = -8.674500000000011
$ws.Range("B27").Value = 5.987300000000009
$ws.Range("D31").Value = -8.470699999999997
$ws.Range("B32").Value = 6.350600000000001
$ws.Range("D35").Value = -8.024900000000001
$ws.Range("B36").Value = 8.977100000000004
$ws.Range("D37").Value = -7.625199999999999
$ws.Range("B38").Value = 5.091900000000001
$ws.Range("D45").Value = -7.684499999999999
$ws.Range("B46").Value = 6.264700000000002
$ws.Range("D52").Value = -7.886399999999997
$ws.Range("B54").Value = 4.894799999999996
$ws.Range("B55").Value = 5.384599999999997
$ws.Range("B56").Value = 5.384099999999998
$ws.Range("D57").Value = -8.468200000000001
$ws.Range("B67").Value = 5.374399999999998
$ws.Range("B69").Value = 5.141099999999997
$ws.Range("B72").Value = 5.137000000000005
$ws.Range("D81").Value = -6.908699999999993
$ws.Range("B83").Value = 5.211699999999999
$ws.Range("D83").Value = -9.0631
$ws.Range("B86").Value = 5.027200000000002
$ws.Range("B91").Value = 5.479600000000001
$ws.Range("B93").Value = 5.517200000000003
$ws.Range("B99").Value = 4.630699999999996
$ws.Range("D100").Value = -8.0068
$ws.Range("D102").Value = -7.743999999999999
